# Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# in column G ("Recorded By"), but only for the specific rows identified
# in the source diff (the final row of each contiguous block of the
# original value is intentionally left untouched, matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Contiguous row ranges (inclusive) in column G that must be updated.
$ranges = @(
    @(2, 7),
    @(16, 17),
    @(22, 23),
    @(37, 38),
    @(43, 44),
    @(58, 59),
    @(64, 65),
    @(79, 80),
    @(85, 90),
    @(99, 100),
    @(105, 110),
    @(119, 120),
    @(125, 130),
    @(139, 140),
    @(145, 150),
    @(159, 160),
    @(165, 170),
    @(179, 180),
    @(185, 186),
    @(200, 201),
    @(206, 207),
    @(221, 222),
    @(227, 228),
    @(242, 243)
)

foreach ($r in $ranges) {
    $startRow = $r[0]
    $endRow = $r[1]
    for ($row = $startRow; $row -le $endRow; $row++) {
        $cell = $ws.Cells.Item($row, 7)
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}

$wb.Save()
